# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" everywhere it
#   appears (Overview!E2:F2, zh-cn!C2, de-de!C2).
# - The status column(s) on each sheet are narrowed to re-fit the new
#   (shorter) text, same as the column-width refresh that a real report
#   regeneration would trigger.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Cell content -----------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value     = $newStatus
$dede.Range("C2").Value     = $newStatus

# --- Column widths (status columns re-fit to the shorter text) --------
$newColumnWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth
$zhcn.Columns.Item(3).ColumnWidth     = $newColumnWidth
$dede.Columns.Item(3).ColumnWidth     = $newColumnWidth
